$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): extend with two new columns P1, Q1 ---
# Copy the formatting (bold/border/centered style) from O1 into the new cells,
# then overwrite their values.
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("P1").Value = 14
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("Q1").Value = 15

# --- Data rows 2-25 ---
# Swap values between columns I/K and K/M (and O mirrors K's old value),
# i.e. I: 1 -> 2, K: 2 -> 1, M: 1 -> 2, O: 2 -> 1, and add new columns P, Q = 2.
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1
$ws.Range("P2:P25").Value = 2
$ws.Range("Q2:Q25").Value = 2
